# Auto-generated Excel COM-interop script applying the Unicorn_Profits leve-profit refresh
# (scheduled runner update across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 156.2
$ws.Range("I4").Value = 82.75
$ws.Range("J4").Value = 450
$ws.Range("K4").Value = 82.75
$ws.Range("L4").Value = 450
$ws.Range("M4").Value = 31.25
$ws.Range("N4").Value = -678

$ws.Range("H53").Value = 328.54544
$ws.Range("I53").Value = 408
$ws.Range("J53").Value = 283.14285
$ws.Range("K53").Value = 408
$ws.Range("L53").Value = 283.14285
$ws.Range("M53").Value = 229
$ws.Range("N53").Value = -1557.14285

$ws.Range("H76").Value = 7647.129
$ws.Range("I76").Value = 8808.944
$ws.Range("K76").Value = 8808.944
$ws.Range("M76").Value = -8493.944

$ws.Range("H79").Value = 7647.129
$ws.Range("I79").Value = 8808.944
$ws.Range("K79").Value = 8808.944
$ws.Range("M79").Value = -7716.944

$ws.Range("H99").Value = 609.8333
$ws.Range("I99").Value = 141
$ws.Range("J99").Value = 1078.6666
$ws.Range("K99").Value = 423
$ws.Range("L99").Value = 3235.9998
$ws.Range("M99").Value = 1075
$ws.Range("N99").Value = -6231.9998


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 997.3333
$ws.Range("I4").Value = 1000
$ws.Range("J4").Value = 996
$ws.Range("K4").Value = 1000
$ws.Range("L4").Value = 996
$ws.Range("M4").Value = -884
$ws.Range("N4").Value = -1228

$ws.Range("H61").Value = 447219.2
$ws.Range("I61").Value = 372310.88
$ws.Range("J61").Value = 559581.7
$ws.Range("K61").Value = 372310.88
$ws.Range("L61").Value = 559581.7
$ws.Range("M61").Value = -372098.88
$ws.Range("N61").Value = -560005.7

$ws.Range("H97").Value = 1237.2333
$ws.Range("I97").Value = 1237.88
$ws.Range("J97").Value = 1234
$ws.Range("K97").Value = 1237.88
$ws.Range("L97").Value = 1234
$ws.Range("M97").Value = -741.8800000000001
$ws.Range("N97").Value = -2226

$ws.Range("H110").Value = 1034.2778
$ws.Range("I110").Value = 1077.4706
$ws.Range("J110").Value = 300
$ws.Range("K110").Value = 1077.4706
$ws.Range("L110").Value = 300
$ws.Range("M110").Value = 967.5293999999999
$ws.Range("N110").Value = -4390

$ws.Range("H136").Value = 447219.2
$ws.Range("I136").Value = 372310.88
$ws.Range("J136").Value = 559581.7
$ws.Range("K136").Value = 1116932.64
$ws.Range("L136").Value = 1678745.1
$ws.Range("M136").Value = -1114382.64
$ws.Range("N136").Value = -1683845.1


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2418.0942
$ws.Range("I105").Value = 2219.9524
$ws.Range("K105").Value = 2219.9524
$ws.Range("M105").Value = -472.9524000000001

$ws.Range("H107").Value = 2101.55
$ws.Range("I107").Value = 2255.4
$ws.Range("J107").Value = 1640
$ws.Range("K107").Value = 2255.4
$ws.Range("L107").Value = 1640
$ws.Range("M107").Value = -335.4000000000001
$ws.Range("N107").Value = -5480

$ws.Range("H134").Value = 52909.43
$ws.Range("I134").Value = 78830.62
$ws.Range("J134").Value = 10787.5
$ws.Range("K134").Value = 236491.86
$ws.Range("L134").Value = 32362.5
$ws.Range("M134").Value = -233956.86
$ws.Range("N134").Value = -37432.5


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 125802.875
$ws.Range("I16").Value = 167435
$ws.Range("K16").Value = 167435
$ws.Range("M16").Value = -167148

$ws.Range("H64").Value = 28900
$ws.Range("J64").Value = 28900
$ws.Range("L64").Value = 28900
$ws.Range("N64").Value = -29396

$ws.Range("H67").Value = 28900
$ws.Range("J67").Value = 28900
$ws.Range("L67").Value = 28900
$ws.Range("N67").Value = -30616

$ws.Range("H107").Value = 378.26086
$ws.Range("I107").Value = 296.75
$ws.Range("J107").Value = 467.18182
$ws.Range("K107").Value = 296.75
$ws.Range("L107").Value = 467.18182
$ws.Range("M107").Value = 1623.25
$ws.Range("N107").Value = -4307.18182

$ws.Range("H113").Value = 125802.875
$ws.Range("I113").Value = 167435
$ws.Range("K113").Value = 167435
$ws.Range("M113").Value = -165265

$ws.Range("H132").Value = 21742422
$ws.Range("I132").Value = 50001292
$ws.Range("K132").Value = 150003876
$ws.Range("M132").Value = -150001346


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 496.12
$ws.Range("J23").Value = 536
$ws.Range("L23").Value = 1608
$ws.Range("N23").Value = -2078

$ws.Range("H34").Value = 1473.9354
$ws.Range("I34").Value = 280
$ws.Range("J34").Value = 1650.8148
$ws.Range("K34").Value = 840
$ws.Range("L34").Value = 4952.4444
$ws.Range("M34").Value = -756
$ws.Range("N34").Value = -5120.4444

$ws.Range("H39").Value = 5510.3706
$ws.Range("I39").Value = 1000
$ws.Range("J39").Value = 6074.1665
$ws.Range("K39").Value = 3000
$ws.Range("L39").Value = 18222.4995
$ws.Range("M39").Value = -2706
$ws.Range("N39").Value = -18810.4995

$ws.Range("H41").Value = 650
$ws.Range("J41").Value = 1000
$ws.Range("L41").Value = 3000
$ws.Range("N41").Value = -3676

$ws.Range("H51").Value = 340
$ws.Range("I51").Value = 340
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 1020
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = -560
$ws.Range("N51").ClearContents() | Out-Null

$ws.Range("H55").Value = 3390
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 3390
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 10170
$ws.Range("M55").ClearContents() | Out-Null
$ws.Range("N55").Value = -10524

$ws.Range("H60").Value = 307.30768
$ws.Range("I60").Value = 177.77777
$ws.Range("J60").Value = 598.75
$ws.Range("K60").Value = 533.33331
$ws.Range("L60").Value = 1796.25
$ws.Range("M60").Value = -282.33331
$ws.Range("N60").Value = -2298.25

$ws.Range("H70").Value = 3701.7144
$ws.Range("I70").Value = 1087.3334
$ws.Range("J70").Value = 5662.5
$ws.Range("K70").Value = 3262.0002
$ws.Range("L70").Value = 16987.5
$ws.Range("M70").Value = -2947.0002
$ws.Range("N70").Value = -17617.5

$ws.Range("H73").Value = 3701.7144
$ws.Range("I73").Value = 1087.3334
$ws.Range("J73").Value = 5662.5
$ws.Range("K73").Value = 3262.0002
$ws.Range("L73").Value = 16987.5
$ws.Range("M73").Value = -2170.0002
$ws.Range("N73").Value = -19171.5

$ws.Range("H97").Value = 317.4
$ws.Range("I97").Value = 270
$ws.Range("J97").Value = 364.8
$ws.Range("K97").Value = 810
$ws.Range("L97").Value = 1094.4
$ws.Range("M97").Value = -314
$ws.Range("N97").Value = -2086.4

$ws.Range("H117").Value = 2628.261
$ws.Range("I117").Value = 3615.5
$ws.Range("J117").Value = 1551.2727
$ws.Range("K117").Value = 10846.5
$ws.Range("L117").Value = 4653.8181
$ws.Range("M117").Value = -7404.5
$ws.Range("N117").Value = -11537.8181


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4056.2812
$ws.Range("I102").Value = 4444.32
$ws.Range("J102").Value = 2670.4285
$ws.Range("K102").Value = 4444.32
$ws.Range("L102").Value = 2670.4285
$ws.Range("M102").Value = -2822.32
$ws.Range("N102").Value = -5914.4285

$ws.Range("H107").Value = 6076.5884
$ws.Range("I107").Value = 7357.857
$ws.Range("J107").Value = 97.333336
$ws.Range("K107").Value = 7357.857
$ws.Range("L107").Value = 97.333336
$ws.Range("M107").Value = -5437.857
$ws.Range("N107").Value = -3937.333336

$ws.Range("H126").Value = 1760.6857
$ws.Range("I126").Value = 1463.96
$ws.Range("J126").Value = 2502.5
$ws.Range("K126").Value = 4391.88
$ws.Range("L126").Value = 7507.5
$ws.Range("M126").Value = -1921.88
$ws.Range("N126").Value = -12447.5


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1522.619
$ws.Range("I61").Value = 1522.619
$ws.Range("K61").Value = 1522.619
$ws.Range("M61").Value = -1320.619

$ws.Range("H113").Value = 1522.619
$ws.Range("I113").Value = 1522.619
$ws.Range("K113").Value = 1522.619
$ws.Range("M113").Value = 647.3810000000001


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 607.625
$ws.Range("I113").Value = 564.6667
$ws.Range("K113").Value = 1694.0001
$ws.Range("M113").Value = 475.9999

$ws.Range("H126").Value = 1790.8966
$ws.Range("I126").Value = 1653.0834
$ws.Range("J126").Value = 2452.4
$ws.Range("K126").Value = 4959.2502
$ws.Range("L126").Value = 7357.200000000001
$ws.Range("M126").Value = -2489.2502
$ws.Range("N126").Value = -12297.2

$ws.Range("H136").Value = 213336.16
$ws.Range("I136").Value = 29794
$ws.Range("J136").Value = 591217.06
$ws.Range("K136").Value = 89382
$ws.Range("L136").Value = 1773651.18
$ws.Range("M136").Value = -86832
$ws.Range("N136").Value = -1778751.18

